# Update 'F' column ('想去人数' / interest counts) across sheets as
# regenerated by the gh-pages data refresh (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 637
$ws.Range("F3").Value = 696
$ws.Range("F4").Value = 942
$ws.Range("F5").Value = 713
$ws.Range("F6").Value = 836
$ws.Range("F7").Value = 398
$ws.Range("F8").Value = 596
$ws.Range("F9").Value = 127
$ws.Range("F10").Value = 1202
$ws.Range("F11").Value = 629
$ws.Range("F12").Value = 381
$ws.Range("F13").Value = 505
$ws.Range("F14").Value = 164
$ws.Range("F16").Value = 444
$ws.Range("F17").Value = 348
$ws.Range("F19").Value = 82
$ws.Range("F20").Value = 553
$ws.Range("F21").Value = 76
$ws.Range("F22").Value = 570
$ws.Range("F24").Value = 740

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 180
$ws.Range("F9").Value = 220
$ws.Range("F13").Value = 90

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 637
$ws.Range("F7").Value = 696
$ws.Range("F8").Value = 942
$ws.Range("F9").Value = 713
$ws.Range("F10").Value = 836
$ws.Range("F11").Value = 398
$ws.Range("F12").Value = 596
$ws.Range("F13").Value = 127
$ws.Range("F14").Value = 1202
$ws.Range("F15").Value = 629
$ws.Range("F18").Value = 381
$ws.Range("F19").Value = 505
$ws.Range("F21").Value = 164
$ws.Range("F23").Value = 444
$ws.Range("F24").Value = 180
$ws.Range("F25").Value = 348
$ws.Range("F27").Value = 82
$ws.Range("F28").Value = 220
$ws.Range("F30").Value = 553
$ws.Range("F33").Value = 90
$ws.Range("F34").Value = 76
$ws.Range("F35").Value = 570
$ws.Range("F37").Value = 740

